# Sniper rounds damage nerf
# Lower the "IRL Speed" (H column) multiplier for the sniper-caliber ammo
# rows on Feuil1. K/E/J columns are formula-driven off H (and I/D which are
# unchanged), so Excel recalculates them automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("H2").Value  = 1.26    # ammo_magnum_300
$ws.Range("H4").Value  = 1       # ammo_7.92x33_fmj
$ws.Range("H5").Value  = 1       # ammo_7.92x33_ap
$ws.Range("H6").Value  = 0.85    # ammo_7.62x54_7h1
$ws.Range("H7").Value  = 0.85    # ammo_7.62x54_ap
$ws.Range("H8").Value  = 0.85    # ammo_7.62x54_7h14
$ws.Range("H9").Value  = 0.89    # ammo_7.62x51_fmj
$ws.Range("H10").Value = 0.89    # ammo_7.62x51_ap
$ws.Range("H11").Value = 1.2     # ammo_12.7x55_fmj
$ws.Range("H12").Value = 1.2     # ammo_12.7x55_ap

$excel.CalculateFull()

# Restore the active cell/selection noted in the diff for this sheet.
$ws.Range("H8").Select()
